$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("F3").Value2 = 1290
$ws.Range("F5").Value2 = 2070
$ws.Range("F6").Value2 = 130
$ws.Range("F7").Value2 = 861
$ws.Range("F10").Value2 = 143
$ws.Range("F11").Value2 = 1072
$ws.Range("F13").Value2 = 36
$ws.Range("F14").Value2 = 680
$ws.Range("F15").Value2 = 1333
$ws.Range("F16").Value2 = 1027
$ws.Range("F19").Value2 = 733
$ws.Range("F20").Value2 = 83
$ws.Range("F24").Value2 = 1225
$ws.Range("F26").Value2 = 431
$ws.Range("F28").Value2 = 5170
$ws.Range("F29").Value2 = 254
$ws.Range("F31").Value2 = 2424
$ws.Range("F32").Value2 = 5814
$ws.Range("F34").Value2 = 976
$ws.Range("F36").Value2 = 59
$ws.Range("F37").Value2 = 1043
$ws.Range("F39").Value2 = 26
$ws.Range("F41").Value2 = 670
$ws.Range("F43").Value2 = 36
$ws.Range("F46").Value2 = 25
$ws.Range("F47").Value2 = 92

$ws = $wb.Worksheets.Item(2)
$ws.Range("F7").Value2 = 11
$ws.Range("F9").Value2 = 464
$ws.Range("F11").Value2 = 97
$ws.Range("F12").Value2 = 91
$ws.Range("F13").Value2 = 124
$ws.Range("F15").Value2 = 660
$ws.Range("F32").Value2 = 147
$ws.Range("F33").Value2 = 100
$ws.Range("F36").Value2 = 65
$ws.Range("F40").Value2 = 892
$ws.Range("F41").Value2 = 486

$ws = $wb.Worksheets.Item(3)
$ws.Range("F4").Value2 = 659
$ws.Range("F5").Value2 = 749
$ws.Range("F6").Value2 = 369
$ws.Range("F7").Value2 = 210

$ws = $wb.Worksheets.Item(4)
$ws.Range("F3").Value2 = 659
$ws.Range("F5").Value2 = 1290
$ws.Range("F7").Value2 = 369
$ws.Range("F8").Value2 = 210
$ws.Range("F9").Value2 = 210
$ws.Range("F10").Value2 = 2070
$ws.Range("F12").Value2 = 861
$ws.Range("F16").Value2 = 143
$ws.Range("F17").Value2 = 1072
$ws.Range("F19").Value2 = 36
$ws.Range("F20").Value2 = 464
$ws.Range("F21").Value2 = 680
$ws.Range("F22").Value2 = 1333
$ws.Range("F23").Value2 = 97
$ws.Range("F24").Value2 = 1027
$ws.Range("F26").Value2 = 124
$ws.Range("F27").Value2 = 733
$ws.Range("F28").Value2 = 83
$ws.Range("F31").Value2 = 1225
$ws.Range("F33").Value2 = 431
$ws.Range("F35").Value2 = 5170
$ws.Range("F36").Value2 = 254
$ws.Range("F38").Value2 = 2424
$ws.Range("F39").Value2 = 5814
$ws.Range("F40").Value2 = 976
$ws.Range("F43").Value2 = 59
$ws.Range("F44").Value2 = 1043
$ws.Range("F45").Value2 = 26
$ws.Range("F46").Value2 = 670
$ws.Range("F47").Value2 = 65
$ws.Range("F48").Value2 = 36
$ws.Range("F49").Value2 = 892
$ws.Range("F50").Value2 = 486
